# Decrement the "Days_Left" value (column H) by 1 for each data row (rows 2-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 8)  # Column H = 8
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current - 1
    }
}
